$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for every touched cell so numeric-looking / percent-looking
# strings (e.g. "332.33", "1.07%") are kept as literal text, matching the source data.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '332.33'
$ws.Range("E2").Value = '1.07%'
$ws.Range("D3").Value = '45.63'
$ws.Range("E3").Value = '2.89%'
$ws.Range("D4").Value = '5.481'
$ws.Range("E4").Value = '-0.23%'
$ws.Range("D5").Value = '0.08527'
$ws.Range("E5").Value = '5.40%'
$ws.Range("D6").Value = '2.066'
$ws.Range("E6").Value = '1.51%'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '0.9885'
$ws.Range("E7").Value = '3.58%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '2.534'
$ws.Range("E8").Value = '-1.78%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.1164'
$ws.Range("E9").Value = '2.13%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1919'
$ws.Range("E10").Value = '1.68%'
$ws.Range("B11").Value = 'MCDex'
$ws.Range("C11").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D11").Value = '9.448'
$ws.Range("E11").Value = '-7.11%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.09793'
$ws.Range("E12").Value = '-1.94%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.04711'
$ws.Range("E13").Value = '-3.27%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.1060'
$ws.Range("E14").Value = '-0.43%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001292'
$ws.Range("E15").Value = '1.27%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005915'
$ws.Range("E16").Value = '2.97%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.392'
$ws.Range("E17").Value = '0.62%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '4.435'
$ws.Range("E18").Value = '0.76%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3356'
$ws.Range("E19").Value = '-1.45%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").Value = '0.1352'
$ws.Range("E20").Value = '-3.34%'
$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D21").Value = '0.2546'
$ws.Range("E21").Value = '-1.30%'
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").Value = '0.04138'
$ws.Range("E22").Value = '1.42%'
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").Value = '0.001300'
$ws.Range("E23").Value = '-0.44%'
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D24").Value = '0.004463'
$ws.Range("E24").Value = '2.32%'
$ws.Range("D25").Value = '0.0001304'
$ws.Range("E25").Value = '4.17%'
$ws.Range("E26").Value = '-20.22%'
$ws.Range("D38").Value = '0.02769'
$ws.Range("E38").Value = '6.76%'
$ws.Range("D39").Value = '0.05741'
$ws.Range("E39").Value = '0.72%'
$ws.Range("D40").Value = '0.007860'
$ws.Range("E40").Value = '3.77%'
$ws.Range("E41").Value = '2.17%'
$ws.Range("D42").Value = '0.007242'
$ws.Range("E42").Value = '-1.37%'
$ws.Range("D43").Value = '0.002121'
$ws.Range("E43").Value = '5.53%'
$ws.Range("D44").Value = '0.007929'
$ws.Range("E44").Value = '-12.78%'
$ws.Range("D45").Value = '0.3407'
$ws.Range("D46").Value = '0.00006983'
$ws.Range("E46").Value = '-0.63%'
$ws.Range("E47").Value = '0.17%'
$ws.Range("E48").Value = '0.37%'
$ws.Range("D49").Value = '0.003459'
$ws.Range("E49").Value = '-1.23%'
$ws.Range("D50").Value = '0.003530'
$ws.Range("E50").Value = '0.80%'
$ws.Range("D51").Value = '0.00002106'
$ws.Range("E51").Value = '0.17%'

# Restore default styling (drop the temporary text-format override) now that
# every value has been committed as text.
$editRange.Style = "Normal"
